# "cost precalculated in scenario"
#
# The id_building_type column (column B) is no longer needed in this
# labor-cost scenario input sheet, so it is removed entirely: every
# subsequent column (id_cooling_technology, unit, and all the year
# columns) shifts one place to the left, and the id_region / the
# (now relocated) id_cooling_technology sample values are updated to
# reflect the re-precalculated scenario row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet hosts a single Excel Table ("Table1") bound to A1:AS2.
# Deleting a column that is part of a Table doesn't reliably resize the
# Table's own definition in this host, so unlist it first, perform the
# plain worksheet column delete (which correctly shifts cells / column
# widths / shared strings), then re-create the Table over the new,
# narrower range so xl/tables/table1.xml stays in sync.
$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name
$lo.Unlist()

# Remove column B (id_building_type) entirely - everything to its right
# (id_cooling_technology, unit, 2010..2050) shifts left by one column.
$ws.Columns("B").Delete()

# Update the sample data row for the now-shifted columns:
# id_region goes from 1 -> 9, and id_cooling_technology (now column B,
# previously column C) goes from 11 -> 1.
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 1

# Re-create the table over the new used range (A1:AR2) with the same name
# and table style as before.
$newlo = $ws.ListObjects.Add(1, $ws.Range("A1:AR2"), [System.Reflection.Missing]::Value, 1)
$newlo.Name = $tableName
$newlo.TableStyle = "TableStyleMedium6"

# Restore the active cell to where the author left off editing.
$null = $ws.Range("C8").Select()
